$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")
$ws.Range("A1").Value = "TEST"
